$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1383.3334
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 1125
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 10125
$ws.Range("M125").Value = -14640
$ws.Range("N125").Value = -15045
$ws.Range("H127").Value = 1019.1429
$ws.Range("I127").Value = 681.6667
$ws.Range("J127").Value = 1272.25
$ws.Range("K127").Value = 2045.0001
$ws.Range("L127").Value = 3816.75
$ws.Range("M127").Value = 2914.9999
$ws.Range("N127").Value = -13736.75
$ws.Range("H132").Value = 20838268
$ws.Range("I132").Value = 23814450
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 71443350
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -71440820
$ws.Range("N132").Value = -20060
$ws.Range("H137").Value = 1664.5667
$ws.Range("I137").Value = 915.93335
$ws.Range("J137").Value = 2413.2
$ws.Range("K137").Value = 2747.80005
$ws.Range("L137").Value = 7239.599999999999
$ws.Range("M137").Value = -197.8000499999998
$ws.Range("N137").Value = -12339.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2607.768
$ws.Range("I32").Value = 2825.6726
$ws.Range("J32").Value = 1751.7142
$ws.Range("K32").Value = 2825.6726
$ws.Range("L32").Value = 1751.7142
$ws.Range("M32").Value = -2538.6726
$ws.Range("N32").Value = -2325.7142
$ws.Range("H45").Value = 1871.2273
$ws.Range("I45").Value = 1719.2354
$ws.Range("J45").Value = 2388
$ws.Range("K45").Value = 1719.2354
$ws.Range("L45").Value = 2388
$ws.Range("M45").Value = -1342.2354
$ws.Range("N45").Value = -3142
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -788
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1391.6666
$ws.Range("I74").Value = 1441.875
$ws.Range("J74").Value = 990
$ws.Range("K74").Value = 1441.875
$ws.Range("L74").Value = 990
$ws.Range("M74").Value = -567.875
$ws.Range("N74").Value = -2738
$ws.Range("H77").Value = 1391.6666
$ws.Range("I77").Value = 1441.875
$ws.Range("J77").Value = 990
$ws.Range("K77").Value = 7209.375
$ws.Range("L77").Value = 4950
$ws.Range("M77").Value = -2841.375
$ws.Range("N77").Value = -13686
$ws.Range("H130").Value = 15429
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 15429
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 15429
$ws.Range("N130").Value = -25469
$ws.Range("H132").Value = 2904.9688
$ws.Range("I132").Value = 2659.5
$ws.Range("J132").Value = 3968.6667
$ws.Range("K132").Value = 7978.5
$ws.Range("L132").Value = 11906.0001
$ws.Range("M132").Value = -5448.5
$ws.Range("N132").Value = -16966.0001
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -14100

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 39116
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 39116
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 39116
$ws.Range("N60").Value = -40314
$ws.Range("H134").Value = 10641.869
$ws.Range("I134").Value = 7660.6113
$ws.Range("J134").Value = 21374.4
$ws.Range("K134").Value = 22981.8339
$ws.Range("L134").Value = 64123.2
$ws.Range("M134").Value = -20446.8339
$ws.Range("N134").Value = -69193.20000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 166667920
$ws.Range("I16").Value = 250001120
$ws.Range("J16").Value = 1496
$ws.Range("K16").Value = 250001120
$ws.Range("L16").Value = 1496
$ws.Range("M16").Value = -250000833
$ws.Range("N16").Value = -2070
$ws.Range("H58").Value = 1680
$ws.Range("I58").Value = 1573.3334
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1573.3334
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1370.3334
$ws.Range("N58").Value = -2406
$ws.Range("H107").Value = 576.1667
$ws.Range("I107").Value = 273.55554
$ws.Range("J107").Value = 757.73334
$ws.Range("K107").Value = 273.55554
$ws.Range("L107").Value = 757.73334
$ws.Range("M107").Value = 1646.44446
$ws.Range("N107").Value = -4597.73334
$ws.Range("H113").Value = 166667920
$ws.Range("I113").Value = 250001120
$ws.Range("J113").Value = 1496
$ws.Range("K113").Value = 250001120
$ws.Range("L113").Value = 1496
$ws.Range("M113").Value = -249998950
$ws.Range("N113").Value = -5836
$ws.Range("H132").Value = 6779.0835
$ws.Range("I132").Value = 7816.1665
$ws.Range("J132").Value = 3667.8333
$ws.Range("K132").Value = 23448.4995
$ws.Range("L132").Value = 11003.4999
$ws.Range("M132").Value = -20918.4995
$ws.Range("N132").Value = -16063.4999
$ws.Range("H134").Value = 11906254
$ws.Range("I134").Value = 14494197
$ws.Range("J134").Value = 1718.8
$ws.Range("K134").Value = 43482591
$ws.Range("L134").Value = 5156.4
$ws.Range("M134").Value = -43480056
$ws.Range("N134").Value = -10226.4
$ws.Range("H136").Value = 1680
$ws.Range("I136").Value = 1573.3334
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4720.0002
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2170.0002
$ws.Range("N136").Value = -11100

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1553.2916
$ws.Range("I34").Value = 650
$ws.Range("J34").Value = 1733.95
$ws.Range("K34").Value = 1950
$ws.Range("L34").Value = 5201.85
$ws.Range("M34").Value = -1866
$ws.Range("N34").Value = -5369.85
$ws.Range("H39").Value = 1750
$ws.Range("I39").Value = 50
$ws.Range("J39").Value = 1992.8572
$ws.Range("K39").Value = 150
$ws.Range("L39").Value = 5978.571599999999
$ws.Range("M39").Value = 144
$ws.Range("N39").Value = -6566.571599999999
$ws.Range("H68").Value = 1862.5
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 2055
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 6165
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -7787
$ws.Range("H71").Value = 1862.5
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 2055
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 18495
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -26607
$ws.Range("H75").Value = 2504.2856
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2504.2856
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 7512.8568
$ws.Range("N75").Value = -9508.856800000001
$ws.Range("H78").Value = 2504.2856
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2504.2856
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 22538.5704
$ws.Range("N78").Value = -32522.5704
$ws.Range("H99").Value = 1900.9231
$ws.Range("I99").Value = 575
$ws.Range("J99").Value = 2490.2222
$ws.Range("K99").Value = 1725
$ws.Range("L99").Value = 7470.6666
$ws.Range("M99").Value = 521
$ws.Range("N99").Value = -11962.6666
$ws.Range("H100").Value = 3348
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3348
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 10044
$ws.Range("N100").Value = -11666
$ws.Range("H131").Value = 20834740
$ws.Range("I131").Value = 100000610
$ws.Range("J131").Value = 1617.0526
$ws.Range("K131").Value = 300001830
$ws.Range("L131").Value = 4851.1578
$ws.Range("M131").Value = -299996790
$ws.Range("N131").Value = -14931.1578
$ws.Range("H139").Value = 1648.9445
$ws.Range("I139").Value = 1426.1364
$ws.Range("J139").Value = 1999.0714
$ws.Range("K139").Value = 4278.4092
$ws.Range("L139").Value = 5997.2142
$ws.Range("M139").Value = 861.5907999999999
$ws.Range("N139").Value = -16277.2142

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 45004810
$ws.Range("I70").Value = 83336600
$ws.Range("J70").Value = 28576900
$ws.Range("K70").Value = 83336600
$ws.Range("L70").Value = 28576900
$ws.Range("M70").Value = -83336330
$ws.Range("N70").Value = -28577440
$ws.Range("H73").Value = 45004810
$ws.Range("I73").Value = 83336600
$ws.Range("J73").Value = 28576900
$ws.Range("K73").Value = 83336600
$ws.Range("L73").Value = 28576900
$ws.Range("M73").Value = -83335664
$ws.Range("N73").Value = -28578772
$ws.Range("H126").Value = 2403.3225
$ws.Range("I126").Value = 1492
$ws.Range("J126").Value = 3846.25
$ws.Range("K126").Value = 4476
$ws.Range("L126").Value = 11538.75
$ws.Range("M126").Value = -2006
$ws.Range("N126").Value = -16478.75
$ws.Range("H132").Value = 3157.375
$ws.Range("I132").Value = 2767
$ws.Range("J132").Value = 3808
$ws.Range("K132").Value = 8301
$ws.Range("L132").Value = 11424
$ws.Range("M132").Value = -5771
$ws.Range("N132").Value = -16484

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1540.1
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1857.2858
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1857.2858
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -2447.2858
$ws.Range("H27").Value = 1540.1
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1857.2858
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1857.2858
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -2071.2858
$ws.Range("H69").Value = 49966.668
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 49966.668
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 49966.668
$ws.Range("N69").Value = -51588.668
$ws.Range("H72").Value = 49966.668
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 49966.668
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 149900.004
$ws.Range("N72").Value = -158012.004
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H136").Value = 5150.5835
$ws.Range("I136").Value = 7513.3335
$ws.Range("J136").Value = 1212.6666
$ws.Range("K136").Value = 22540.0005
$ws.Range("L136").Value = 3637.9998
$ws.Range("M136").Value = -19990.0005
$ws.Range("N136").Value = -8737.9998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 7495
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7495
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 7495
$ws.Range("N93").Value = -12487
$ws.Range("H112").Value = 28000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 28000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 28000
$ws.Range("N112").Value = -30954
$ws.Range("H122").Value = 28892844
$ws.Range("I122").Value = 28892844
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 86678532
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -86676082
$ws.Range("H123").Value = 58502.8
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 58502.8
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 58502.8
$ws.Range("N123").Value = -68302.8
$ws.Range("H126").Value = 48309864
$ws.Range("I126").Value = 52910664
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 158731992
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -158729522
$ws.Range("N126").Value = -9290
$ws.Range("H132").Value = 3939.3333
$ws.Range("I132").Value = 4243.5454
$ws.Range("J132").Value = 2600.8
$ws.Range("K132").Value = 12730.6362
$ws.Range("L132").Value = 7802.400000000001
$ws.Range("M132").Value = -10200.6362
$ws.Range("N132").Value = -12862.4
$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 750
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2250
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 300
$ws.Range("N136").Value = -14100
$ws.Range("H138").Value = 50000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60284
